$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B2").Value = "var001_Login_001_Successful"
$ws.Range("B3").Value = "var001_Login_002_LoginWithoutEmailAndPassword"
$ws.Range("B4").Value = "var002_PasswordReset_001_Successful"
$ws.Range("B5").Value = "var002_PasswordReset_002_BackToLogin"

$ws.Range("B10").Select()
